# Guest Authentication Through Our Partners - communication plan update
# Applies the edits described in the commit "update communication plan for partners"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Last Updated date: 9/12/24 -> 9/13/24
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(2)
$p.Range.Find.Execute("Last Updated:  9/12/24", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "Last Updated:  9/13/24", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. "Note: decision will be made..." -> "Assumption:  " intro, and the
#    trailing sentence about Okta/Auth0 -> Auth0 confirmation wording.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(8)
$p.Range.Find.Execute("Note: decision will be made on 9/13/24 whether ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "Assumption:  ", 2) | Out-Null
$p = $d.Paragraphs.Item(8)
$p.Range.Find.Execute("will use Okta or Auth0 for the underlying Identity platform.  Assuming we move forward with Auth0:", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, "will be moving forward to Auth0 platform as the underlying Identity platform:", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Merge the "to your developer or integration team" / "to determine
#    level of effort..." runs (no text change, touch to normalize).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(10)
$p.Range.Find.Execute(" to your developer or integration team to determine level of effort for the Partner to update SSO configuration from Ping to Auth0.", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, " to your developer or integration team to determine level of effort for the Partner to update SSO configuration from Ping to Auth0.", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Merge "Create the Partner application in " / "our Identity
#    Provider" runs (no text change, touch to normalize).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(13)
$p.Range.Find.Execute("Create the Partner application in our Identity Provider", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, "Create the Partner application in our Identity Provider", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Insert a blank paragraph right after "See README.md..." and before
#    the "Timeline" heading.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(46)
$p.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 6. Remove the "[Neha, please review and/or revise dates]" placeholder
#    paragraph entirely (now shifted down by 1 due to the blank
#    paragraph inserted above).
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(49)
$p.Range.Delete()

# ---------------------------------------------------------------------
# 7. Push every Timeline date forward by one week and fold the
#    standalone "PROD" runs into the surrounding sentences.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(49)
$p.Range.Find.Execute("9/20/24 ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "9/27/24 ", 2) | Out-Null

$p = $d.Paragraphs.Item(50)
$p.Range.Find.Execute("9/27/24 – Partner communicates with Neha on ETA for when they will be able to update TEST configuration from PingFederate to ", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, "10/4/24 – Partner communicates with Neha on ETA for when they will be able to update TEST configuration from PingFederate to ", 2) | Out-Null

$p = $d.Paragraphs.Item(51)
$p.Range.Find.Execute("10/11/24 – ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "10/18/24 – ", 2) | Out-Null
$p = $d.Paragraphs.Item(51)
$p.Range.Find.Execute(" creates applications in PROD environment for each Partner and provides client id to Partner", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, " creates applications in PROD environment for each Partner and provides client id to Partner", 2) | Out-Null

$p = $d.Paragraphs.Item(52)
$p.Range.Find.Execute("10/18/24 - Partner communicates with Neha on ETA for when they will be able to update PROD configuration from PingFederate to ", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, "10/25/24 - Partner communicates with Neha on ETA for when they will be able to update PROD configuration from PingFederate to ", 2) | Out-Null
$p = $d.Paragraphs.Item(52)
$p.Range.Find.Execute(" Identity and to begin end to end integration testing in their PROD environment", `
                       $true, $false, $false, $false, $false, `
                       $true, 1, $false, " Identity and to begin end to end integration testing in their PROD environment", 2) | Out-Null

Write-Output "done"
